$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 0.8540296666666668
$ws.Range("H2").Value = 2.562089
$ws.Range("M2").Value = 1.660421
$ws.Range("N2").Value = 4.981262999999999
$ws.Range("O2").Value = 0.03714789785507311
$ws.Range("P2").Value = 0.03714789785507311
$ws.Range("Q2").Value = 1.418048793156333
$ws.Range("R2").Value = 12.762439138407
$ws.Range("S2").Value = 0.03714789785507311
$ws.Range("T2").Value = 0.03714789785507311

# Row 3
$ws.Range("G3").Value = 0.8540296666666668
$ws.Range("H3").Value = 2.562089
$ws.Range("O3").Value = 0.5631392661118858
$ws.Range("P3").Value = 0.5631392661118859
$ws.Range("Q3").Value = 21.49674686315656
$ws.Range("R3").Value = 193.470721768409
$ws.Range("S3").Value = 0.5631392661118858
$ws.Range("T3").Value = 0.5631392661118859

# Row 4
$ws.Range("G4").Value = 0.8540296666666668
$ws.Range("H4").Value = 2.562089
$ws.Range("M4").Value = 17.866195
$ws.Range("N4").Value = 53.598585
$ws.Range("O4").Value = 0.399712836033041
$ws.Range("P4").Value = 0.399712836033041
$ws.Range("Q4").Value = 15.25826056045167
$ws.Range("R4").Value = 137.324345044065
$ws.Range("S4").Value = 0.399712836033041
$ws.Range("T4").Value = 0.399712836033041
